$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "rxxx"
$ws.Range("B6").Value = "kyle"
$ws.Range("C6").Value = "we should update belt capstan pressure to 20psi"
$ws.Range("D6").Value = "2025-09-30 13:09:26"
